$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.377.72"
$ws.Range("E2").Value = "  -2.68%  "
$ws.Range("D3").Value = "1.981.57"
$ws.Range("E3").Value = "  -3.43%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.51"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -8.87%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0821"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.16%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.863"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.78%  "
$ws.Range("E15").Value = "  -5.69%  "
$ws.Range("D16").Value = "2.270.39"
$ws.Range("E16").Value = "  -3.43%  "
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").Value = "1.980.31"
$ws.Range("E18").Value = "  -3.41%  "
$ws.Range("D19").Value = "36.289.52"
$ws.Range("E19").Value = "  -2.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.94%  "
$ws.Range("D21").Value = "0.0₃0866"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("E22").Value = "  -3.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.132"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.79%  "
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0626"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.86%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.80%  "
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0967"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.21%  "
$ws.Range("E44").Value = "  -2.06%  "
$ws.Range("E45").Value = "  -4.78%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.53%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "92.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.12%  "
$ws.Range("D48").Value = "1.371.98"
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("E49").Value = "  -5.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.54%  "
